$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---- Column widths (columns I, J, K) ----
$ws.Columns(9).ColumnWidth = 36.85
$ws.Columns(10).ColumnWidth = 34.0
$ws.Columns(11).ColumnWidth = 34.25

# ---- Apply styles first (copy/paste formats keeps the style table clean, no orphans) ----

# Header row (row 1): copy bold+centered+wrap style from an existing header cell
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial($xlPasteFormats)
$ws.Range("J1").PasteSpecial($xlPasteFormats)
$ws.Range("K1").PasteSpecial($xlPasteFormats)

# Body cells: copy wrap-only style from an existing body wrap cell
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial($xlPasteFormats)
$ws.Range("I3").PasteSpecial($xlPasteFormats)
$ws.Range("I4").PasteSpecial($xlPasteFormats)
$ws.Range("I5").PasteSpecial($xlPasteFormats)
$ws.Range("I6").PasteSpecial($xlPasteFormats)
$ws.Range("I7").PasteSpecial($xlPasteFormats)
$ws.Range("I8").PasteSpecial($xlPasteFormats)
$ws.Range("J2").PasteSpecial($xlPasteFormats)
$ws.Range("K2").PasteSpecial($xlPasteFormats)

# Quote-prefixed body cells: create the new quotePrefix+wrap style once on J3, then replicate it.
$ws.Range("J3").WrapText = $true

$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial($xlPasteFormats)
$ws.Range("K4").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---- Assign cell values (order matches original authoring order for shared strings) ----
$ws.Range("I1").Value = "Max_Habscore"
$ws.Range("J1").Value = "Warning"
$ws.Range("J3").Value = "'Difference +2') Difference between Habscore_a and Habscore_b is >=2; score difference is worth double checking."
$ws.Range("J2").Value = "NA) Difference between Habscore_a and Habscore_b is <2; score difference is not notable."
$ws.Range("K1").Value = "Check_Warnings"
$ws.Range("K2").Value = "NA) Warning column did not produce a warning that recommends a second look."
$ws.Range("K3").Value = "'Good') Warning was double checked for the site and the scores were appropriately assigned."
$ws.Range("K4").Value = "'REVISIT') Warning was double checked for the site and the scores were NOT appropriately assigned."

$ws.Range("I2").Value = "0) No oysters: bare sand or mud."
$ws.Range("I3").Value = "1) Sand or mud with habitat structure (oysters, shell, granite, rock, tunicates*or hard substrate) covering <50% of the bottom."
$ws.Range("I4").Value = "2) Habitat structure (oysters, shell, granite, rock, tunicates*or hard substrate) covering >50% of the bottom, structure much less than the shell height of an individual oyster."
$ws.Range("I5").Value = "3**) Habitat structure (oysters, shell, granite, rock, tunicates*or hard substrate) covering >50% of the bottom and structure height equal to or greater than the shell height of an individual oyster."
$ws.Range("I6").Value = "9) Poor quality, can't assign an accurate habitat score."
$ws.Range("I7").Value = "* Tunicates need hard substrate to grow on, so if the photo is full of tunicates it is assumed the bottom is hard substrate."
$ws.Range("I8").Value = "** Often scores of 3 also have oysters growing in clusters."

# ---- Row heights ----
$ws.Rows(2).RowHeight = 43.5

# ---- View state: scroll so column G is the leftmost visible column, select J4 ----
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$ws.Range("J4").Select()
